$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C9) from 45174 (2023-09-05) to 45175 (2023-09-06)
foreach ($row in 2..9) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
